$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 1372.1666
$ws.Range("I33").Value = 1790.25
$ws.Range("K33").Value = 1790.25
$ws.Range("M33").Value = -1561.25

# Row 127
$ws.Range("H127").Value = 629.7143
$ws.Range("I127").Value = 681.8333
$ws.Range("J127").Value = 317
$ws.Range("K127").Value = 2045.4999
$ws.Range("L127").Value = 951
$ws.Range("M127").Value = 2914.5001
$ws.Range("N127").Value = -10871

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 518.33
$ws.Range("I32").Value = 473.54254
$ws.Range("K32").Value = 473.54254
$ws.Range("M32").Value = -186.54254

# Row 33
$ws.Range("H33").Value = 20323.75
$ws.Range("I33").Value = 15718.2
$ws.Range("J33").Value = 27999.666
$ws.Range("K33").Value = 15718.2
$ws.Range("L33").Value = 27999.666
$ws.Range("M33").Value = -15389.2
$ws.Range("N33").Value = -28657.666

# Row 36
$ws.Range("H36").Value = 8749.75
$ws.Range("I36").Value = 8333
$ws.Range("K36").Value = 8333
$ws.Range("M36").Value = -7987

# Row 38
$ws.Range("H38").Value = 7649.6665
$ws.Range("I38").Value = 3000
$ws.Range("J38").Value = 9974.5
$ws.Range("K38").Value = 3000
$ws.Range("L38").Value = 9974.5
$ws.Range("M38").Value = -2533
$ws.Range("N38").Value = -10908.5

# Row 97
$ws.Range("H97").Value = 53865
$ws.Range("I97").Value = 2862.25
$ws.Range("J97").Value = 257876
$ws.Range("K97").Value = 2862.25
$ws.Range("L97").Value = 257876
$ws.Range("M97").Value = -2366.25
$ws.Range("N97").Value = -258868

# Row 141
$ws.Range("H141").Value = 90280.28999999999
$ws.Range("I141").Value = 45000
$ws.Range("J141").Value = 97827
$ws.Range("K141").Value = 45000
$ws.Range("L141").Value = 97827
$ws.Range("M141").Value = -39820
$ws.Range("N141").Value = -108187

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 2186.9285
$ws.Range("I22").Value = 2068.889
$ws.Range("J22").Value = 2399.4
$ws.Range("K22").Value = 2068.889
$ws.Range("L22").Value = 2399.4
$ws.Range("M22").Value = -1718.889
$ws.Range("N22").Value = -3099.4

# Row 31
$ws.Range("H31").Value = 3406.375
$ws.Range("I31").Value = 2356.25
$ws.Range("J31").Value = 4456.5
$ws.Range("K31").Value = 2356.25
$ws.Range("L31").Value = 4456.5
$ws.Range("M31").Value = -2061.25
$ws.Range("N31").Value = -5046.5

# Row 34
$ws.Range("H34").Value = 3406.375
$ws.Range("I34").Value = 2356.25
$ws.Range("J34").Value = 4456.5
$ws.Range("K34").Value = 2356.25
$ws.Range("L34").Value = 4456.5
$ws.Range("M34").Value = -2154.25
$ws.Range("N34").Value = -4860.5

# Row 52
$ws.Range("H52").Value = 49363.5
$ws.Range("J52").Value = 49363.5
$ws.Range("L52").Value = 49363.5
$ws.Range("N52").Value = -49951.5

# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# Row 105
$ws.Range("H105").Value = 3599.2666
$ws.Range("I105").Value = 1999.8
$ws.Range("J105").Value = 4399
$ws.Range("K105").Value = 1999.8
$ws.Range("L105").Value = 4399
$ws.Range("M105").Value = -252.8
$ws.Range("N105").Value = -7893

# Row 132
$ws.Range("H132").Value = 7564.125
$ws.Range("I132").Value = 3175.75
$ws.Range("J132").Value = 29506
$ws.Range("K132").Value = 9527.25
$ws.Range("L132").Value = 88518
$ws.Range("M132").Value = -6997.25
$ws.Range("N132").Value = -93578

# Row 135
$ws.Range("H135").Value = 101462.125
$ws.Range("J135").Value = 101462.125
$ws.Range("L135").Value = 101462.125
$ws.Range("N135").Value = -111602.125

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 3332.6667
$ws.Range("I68").Value = 2998
$ws.Range("J68").Value = 3399.6
$ws.Range("K68").Value = 8994
$ws.Range("L68").Value = 10198.8
$ws.Range("M68").Value = -8183
$ws.Range("N68").Value = -11820.8

# Row 71
$ws.Range("H71").Value = 3332.6667
$ws.Range("I71").Value = 2998
$ws.Range("J71").Value = 3399.6
$ws.Range("K71").Value = 26982
$ws.Range("L71").Value = 30596.4
$ws.Range("M71").Value = -22926
$ws.Range("N71").Value = -38708.39999999999

# Row 94
$ws.Range("H94").Value = 100002800
$ws.Range("J94").Value = 116668270
$ws.Range("L94").Value = 350004810
$ws.Range("N94").Value = -350006162

# Row 113
$ws.Range("H113").Value = 1381.95
$ws.Range("I113").Value = 952.5
$ws.Range("J113").Value = 1668.25
$ws.Range("K113").Value = 2857.5
$ws.Range("L113").Value = 5004.75
$ws.Range("M113").Value = -687.5
$ws.Range("N113").Value = -9344.75

# Row 128
$ws.Range("H128").Value = 178580
$ws.Range("I128").Value = 178580
$ws.Range("K128").Value = 535740
$ws.Range("M128").Value = -530760

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 9010
$ws.Range("I31").Value = 3030
$ws.Range("J31").Value = 12000
$ws.Range("K31").Value = 3030
$ws.Range("L31").Value = 12000
$ws.Range("M31").Value = -2738
$ws.Range("N31").Value = -12584

# Row 37
$ws.Range("H37").Value = 9010
$ws.Range("I37").Value = 3030
$ws.Range("J37").Value = 12000
$ws.Range("K37").Value = 3030
$ws.Range("L37").Value = 12000
$ws.Range("M37").Value = -2753
$ws.Range("N37").Value = -12554

# Row 135
$ws.Range("H135").Value = 156228.67
$ws.Range("J135").Value = 157043.62
$ws.Range("L135").Value = 157043.62
$ws.Range("N135").Value = -167183.62

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 757.8889
$ws.Range("I22").Value = 728.25
$ws.Range("K22").Value = 728.25
$ws.Range("M22").Value = -433.25

# Row 27
$ws.Range("H27").Value = 757.8889
$ws.Range("I27").Value = 728.25
$ws.Range("K27").Value = 728.25
$ws.Range("M27").Value = -621.25

# Row 32
$ws.Range("H32").Value = 5000
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

# Row 134
$ws.Range("H134").Value = 118798
$ws.Range("J134").Value = 118798
$ws.Range("L134").Value = 118798
$ws.Range("N134").Value = -128938

# Row 136
$ws.Range("H136").Value = 5309986.5
$ws.Range("I136").Value = 7514607
$ws.Range("J136").Value = 18896.3
$ws.Range("K136").Value = 22543821
$ws.Range("L136").Value = 56688.89999999999
$ws.Range("M136").Value = -22541271
$ws.Range("N136").Value = -61788.89999999999

# Row 137
$ws.Range("H137").Value = 139499.25
$ws.Range("J137").Value = 119332.664
$ws.Range("L137").Value = 119332.664
$ws.Range("N137").Value = -129532.664

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 9937.25
$ws.Range("I62").Value = 8200
$ws.Range("K62").Value = 8200
$ws.Range("M62").Value = -7576

# Row 65
$ws.Range("H65").Value = 9937.25
$ws.Range("I65").Value = 8200
$ws.Range("K65").Value = 41000
$ws.Range("M65").Value = -37880

# Row 126
$ws.Range("H126").Value = 10256.833
$ws.Range("I126").Value = 5636.375
$ws.Range("J126").Value = 19497.75
$ws.Range("K126").Value = 16909.125
$ws.Range("L126").Value = 58493.25
$ws.Range("M126").Value = -14439.125
$ws.Range("N126").Value = -63433.25

# Row 133
$ws.Range("H133").Value = 86981.836
$ws.Range("J133").Value = 86981.836
$ws.Range("L133").Value = 86981.836
$ws.Range("N133").Value = -97101.836

# Row 135
$ws.Range("H135").Value = 173072.83
$ws.Range("J135").Value = 173072.83
$ws.Range("L135").Value = 173072.83
$ws.Range("N135").Value = -183212.83
